$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    paragraph (Heading1), matching:
#      <w:p>
#        <w:r/>
#        <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#        <w:r><w:t>: Discover ...</w:t></w:r>
#      </w:p>
# ---------------------------------------------------------------------

$metaRest = ": Discover the thrilling gameplay and haunting graphics of Dark King: Forbidden Riches. Play for free and win up to 2,000x your bet with sticky wilds and free spins."

$titlePara = $d.Paragraphs(1)
$insertPoint = $titlePara.Range.End
$collapsed = $d.Range($insertPoint, $insertPoint)

# Inserting raw WordprocessingML runs via a collapsed range at a paragraph
# boundary creates the new paragraph (no inherited Heading1 style) while an
# extra trailing empty paragraph absorbs the following paragraph's own
# pPr/style untouched; that helper paragraph is removed afterwards.
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$metaXml = '<w:p ' + $ns + '><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>' + $metaRest + '</w:t></w:r></w:p><w:p ' + $ns + '></w:p>'
$collapsed.InsertXML($metaXml)

# Remove the helper empty paragraph left behind between the new paragraph
# and the paragraph that used to follow the title.
$d.Paragraphs(3).Range.Delete()

# ---------------------------------------------------------------------
# 2) Drop the duplicated bold "Play Dark King..." paragraph that used to
#    sit right before the closing italic meta-description paragraph.
# ---------------------------------------------------------------------

$count = $d.Paragraphs.Count
$titleDupe = $d.Paragraphs($count - 1)
$titleDupe.Range.Delete()

# ---------------------------------------------------------------------
# 3) Replace the final (italic) paragraph's text with the new DALL-E
#    image-prompt text, preserving its run formatting / leading empty run.
# ---------------------------------------------------------------------

$newImagePrompt = "DALLE, please create an eye-catching feature image for NetEnt's Dark King: Forbidden Riches game. The image should be in cartoon style and feature a happy Maya warrior with glasses. This character should be a prominent figure in the image, with dark and ominous graphics in the background. The image should capture the fantasy theme of the game and include elements such as skulls, an armored knight, and medieval-style symbols. Please use warm colors that pop and create a sense of excitement and adventure. Overall, the image should entice players to try out this exciting new game."

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastTextRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End - 1)
$lastTextRange.Text = $newImagePrompt

Write-Output "done"
